$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3271.6667
$ws.Range("I32").Value = 2299.25
$ws.Range("K32").Value = 2299.25
$ws.Range("M32").Value = -1973.25
$ws.Range("H42").Value = 2349.8
$ws.Range("I42").Value = 687.25
$ws.Range("J42").Value = 9000
$ws.Range("K42").Value = 2061.75
$ws.Range("L42").Value = 27000
$ws.Range("M42").Value = -1831.75
$ws.Range("N42").Value = -27460
$ws.Range("H64").Value = 3997
$ws.Range("I64").Value = 3997
$ws.Range("K64").Value = 3997
$ws.Range("M64").Value = -3749
$ws.Range("H67").Value = 3997
$ws.Range("I67").Value = 3997
$ws.Range("K67").Value = 3997
$ws.Range("M67").Value = -3139
$ws.Range("H75").Value = 149000
$ws.Range("J75").Value = 149000
$ws.Range("L75").Value = 149000
$ws.Range("N75").Value = -150872
$ws.Range("H78").Value = 149000
$ws.Range("J78").Value = 149000
$ws.Range("L78").Value = 447000
$ws.Range("N78").Value = -456360
$ws.Range("H98").Value = 1979
$ws.Range("I98").Value = 1970.3
$ws.Range("J98").Value = 1993.5
$ws.Range("K98").Value = 1970.3
$ws.Range("L98").Value = 1993.5
$ws.Range("M98").Value = -472.3
$ws.Range("N98").Value = -4989.5
$ws.Range("H112").Value = 4972.5
$ws.Range("I112").Value = 4972
$ws.Range("K112").Value = 14916
$ws.Range("M112").Value = -13808
$ws.Range("H122").Value = 1979
$ws.Range("I122").Value = 1970.3
$ws.Range("J122").Value = 1993.5
$ws.Range("K122").Value = 5910.9
$ws.Range("L122").Value = 5980.5
$ws.Range("M122").Value = -3460.9
$ws.Range("N122").Value = -10880.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2838.6
$ws.Range("I97").Value = 1048.25
$ws.Range("J97").Value = 10000
$ws.Range("K97").Value = 1048.25
$ws.Range("L97").Value = 10000
$ws.Range("M97").Value = -552.25
$ws.Range("N97").Value = -10992
$ws.Range("H122").Value = 4397.875
$ws.Range("I122").Value = 4397.875
$ws.Range("K122").Value = 13193.625
$ws.Range("M122").Value = -10743.625
$ws.Range("H132").Value = 4456.7837
$ws.Range("I132").Value = 3806.3572
$ws.Range("J132").Value = 6480.3335
$ws.Range("K132").Value = 11419.0716
$ws.Range("L132").Value = 19441.0005
$ws.Range("M132").Value = -8889.071599999999
$ws.Range("N132").Value = -24501.0005
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4588.905
$ws.Range("I86").Value = 2276.3076
$ws.Range("K86").Value = 2276.3076
$ws.Range("M86").Value = -1153.3076
$ws.Range("H89").Value = 4588.905
$ws.Range("I89").Value = 2276.3076
$ws.Range("K89").Value = 11381.538
$ws.Range("M89").Value = -5765.538
$ws.Range("H94").Value = 2378.5334
$ws.Range("I94").Value = 1765
$ws.Range("K94").Value = 1765
$ws.Range("M94").Value = -1314
$ws.Range("H99").Value = 1771.125
$ws.Range("I99").Value = 1061.6666
$ws.Range("K99").Value = 1061.6666
$ws.Range("M99").Value = 436.3334
$ws.Range("H105").Value = 3000
$ws.Range("I105").Value = 3000
$ws.Range("K105").Value = 3000
$ws.Range("M105").Value = -1253
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2985.7144
$ws.Range("I31").Value = 2982.6
$ws.Range("K31").Value = 2982.6
$ws.Range("M31").Value = -2687.6
$ws.Range("H34").Value = 2985.7144
$ws.Range("I34").Value = 2982.6
$ws.Range("K34").Value = 2982.6
$ws.Range("M34").Value = -2780.6
$ws.Range("H58").Value = 2267.2727
$ws.Range("I58").Value = 1565.875
$ws.Range("J58").Value = 4137.6665
$ws.Range("K58").Value = 1565.875
$ws.Range("L58").Value = 4137.6665
$ws.Range("M58").Value = -1362.875
$ws.Range("N58").Value = -4543.6665
$ws.Range("H62").Value = 3325.6
$ws.Range("I62").Value = 3325.6
$ws.Range("K62").Value = 3325.6
$ws.Range("M62").Value = -2701.6
$ws.Range("H65").Value = 3325.6
$ws.Range("I65").Value = 3325.6
$ws.Range("K65").Value = 16628
$ws.Range("M65").Value = -13508
$ws.Range("H132").Value = 2782.68
$ws.Range("I132").Value = 2535.35
$ws.Range("J132").Value = 3772
$ws.Range("K132").Value = 7606.049999999999
$ws.Range("L132").Value = 11316
$ws.Range("M132").Value = -5076.049999999999
$ws.Range("N132").Value = -16376
$ws.Range("H136").Value = 2267.2727
$ws.Range("I136").Value = 1565.875
$ws.Range("J136").Value = 4137.6665
$ws.Range("K136").Value = 4697.625
$ws.Range("L136").Value = 12412.9995
$ws.Range("M136").Value = -2147.625
$ws.Range("N136").Value = -17512.9995
$ws.Range("H141").Value = 419999.66
$ws.Range("J141").Value = 584999.5
$ws.Range("L141").Value = 584999.5
$ws.Range("N141").Value = -595359.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 844.5
$ws.Range("I5").Value = 888.75
$ws.Range("J5").Value = 667.5
$ws.Range("K5").Value = 2666.25
$ws.Range("L5").Value = 2002.5
$ws.Range("M5").Value = -2554.25
$ws.Range("N5").Value = -2226.5
$ws.Range("H113").Value = 1900.2307
$ws.Range("I113").Value = 2304.5
$ws.Range("J113").Value = 1553.7142
$ws.Range("K113").Value = 6913.5
$ws.Range("L113").Value = 4661.142599999999
$ws.Range("M113").Value = -4743.5
$ws.Range("N113").Value = -9001.142599999999
$ws.Range("H135").Value = 844.5
$ws.Range("I135").Value = 888.75
$ws.Range("J135").Value = 667.5
$ws.Range("K135").Value = 7998.75
$ws.Range("L135").Value = 6007.5
$ws.Range("M135").Value = -5463.75
$ws.Range("N135").Value = -11077.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12324.5
$ws.Range("J70").Value = 12199.5
$ws.Range("L70").Value = 12199.5
$ws.Range("N70").Value = -12739.5
$ws.Range("H73").Value = 12324.5
$ws.Range("J73").Value = 12199.5
$ws.Range("L73").Value = 12199.5
$ws.Range("N73").Value = -14071.5
$ws.Range("H80").Value = 3999.6667
$ws.Range("I80").Value = 3666
$ws.Range("J80").Value = 4333.3335
$ws.Range("K80").Value = 3666
$ws.Range("L80").Value = 4333.3335
$ws.Range("M80").Value = -2668
$ws.Range("N80").Value = -6329.3335
$ws.Range("H83").Value = 3999.6667
$ws.Range("I83").Value = 3666
$ws.Range("J83").Value = 4333.3335
$ws.Range("K83").Value = 18330
$ws.Range("L83").Value = 21666.6675
$ws.Range("M83").Value = -13338
$ws.Range("N83").Value = -31650.6675
$ws.Range("H97").Value = 1882.2307
$ws.Range("J97").Value = 1798.4
$ws.Range("L97").Value = 1798.4
$ws.Range("N97").Value = -2790.4
$ws.Range("H122").Value = 2204.8333
$ws.Range("I122").Value = 2187.4707
$ws.Range("K122").Value = 6562.4121
$ws.Range("M122").Value = -4112.4121
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 7399.5713
$ws.Range("I68").Value = 7633
$ws.Range("J68").Value = 5999
$ws.Range("K68").Value = 7633
$ws.Range("L68").Value = 5999
$ws.Range("M68").Value = -6884
$ws.Range("N68").Value = -7497
$ws.Range("H71").Value = 7399.5713
$ws.Range("I71").Value = 7633
$ws.Range("J71").Value = 5999
$ws.Range("K71").Value = 38165
$ws.Range("L71").Value = 29995
$ws.Range("M71").Value = -34421
$ws.Range("N71").Value = -37483
$ws.Range("H82").Value = 1349
$ws.Range("I82").Value = 579.1
$ws.Range("J82").Value = 2632.1667
$ws.Range("K82").Value = 579.1
$ws.Range("L82").Value = 2632.1667
$ws.Range("M82").Value = -218.1
$ws.Range("N82").Value = -3354.1667
$ws.Range("H85").Value = 1349
$ws.Range("I85").Value = 579.1
$ws.Range("J85").Value = 2632.1667
$ws.Range("K85").Value = 579.1
$ws.Range("L85").Value = 2632.1667
$ws.Range("M85").Value = 668.9
$ws.Range("N85").Value = -5128.1667
$ws.Range("H100").Value = 1999
$ws.Range("I100").Value = 1999
$ws.Range("K100").Value = 1999
$ws.Range("M100").Value = -1458
$ws.Range("H132").Value = 5450
$ws.Range("I132").Value = 5400
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 16200
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = -13670
$ws.Range("N132").Value = -21560
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 139495
$ws.Range("I29").Value = 139495
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 139495
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -139205
$ws.Range("N29").ClearContents()
$ws.Range("H62").Value = 7103.5
$ws.Range("I62").Value = 6524.4
$ws.Range("J62").Value = 9999
$ws.Range("K62").Value = 6524.4
$ws.Range("L62").Value = 9999
$ws.Range("M62").Value = -5900.4
$ws.Range("N62").Value = -11247
$ws.Range("H65").Value = 7103.5
$ws.Range("I65").Value = 6524.4
$ws.Range("J65").Value = 9999
$ws.Range("K65").Value = 32622
$ws.Range("L65").Value = 49995
$ws.Range("M65").Value = -29502
$ws.Range("N65").Value = -56235
$ws.Range("H81").Value = 2020.9166
$ws.Range("I81").Value = 2020.9166
$ws.Range("K81").Value = 4041.8332
$ws.Range("M81").Value = -2980.8332
$ws.Range("H84").Value = 2020.9166
$ws.Range("I84").Value = 2020.9166
$ws.Range("K84").Value = 20209.166
$ws.Range("M84").Value = -14905.166
$ws.Range("H122").Value = 2091.5625
$ws.Range("I122").Value = 2104.7144
$ws.Range("K122").Value = 6314.1432
$ws.Range("M122").Value = -3864.1432
$ws.Range("H132").Value = 3907.9167
$ws.Range("I132").Value = 3354.0908
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 10062.2724
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -7532.2724
$ws.Range("N132").Value = -35060
$ws.Range("H136").Value = 14002.45
$ws.Range("I136").Value = 10678.75
$ws.Range("J136").Value = 27297.25
$ws.Range("K136").Value = 32036.25
$ws.Range("L136").Value = 81891.75
$ws.Range("M136").Value = -29486.25
$ws.Range("N136").Value = -86991.75
